# Auto-generated edit script: updates computed market-profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) pulled in by the
# scheduled market-data refresh, across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1083.8977
$ws.Range("J129").Value = 1083.8977
$ws.Range("L129").Value = 3251.6931
$ws.Range("N129").Value = -13251.6931
$ws.Range("H135").Value = 596.8511
$ws.Range("I135").Value = 652.3077
$ws.Range("J135").Value = 528.1905
$ws.Range("K135").Value = 5870.7693
$ws.Range("L135").Value = 4753.7145
$ws.Range("M135").Value = -3335.7693
$ws.Range("N135").Value = -9823.7145
$ws.Range("H137").Value = 2615.5134
$ws.Range("I137").Value = 1536.2941
$ws.Range("J137").Value = 3532.85
$ws.Range("K137").Value = 4608.8823
$ws.Range("L137").Value = 10598.55
$ws.Range("M137").Value = -2058.8823
$ws.Range("N137").Value = -15698.55
$ws.Range("H138").Value = 3059.182
$ws.Range("I138").Value = 1994.6
$ws.Range("J138").Value = 3418.838
$ws.Range("K138").Value = 5983.799999999999
$ws.Range("L138").Value = 10256.514
$ws.Range("M138").Value = -843.7999999999993
$ws.Range("N138").Value = -20536.514

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 29345.742
$ws.Range("I2").Value = 798.75
$ws.Range("J2").Value = 143533.72
$ws.Range("K2").Value = 798.75
$ws.Range("L2").Value = 143533.72
$ws.Range("M2").Value = -685.75
$ws.Range("N2").Value = -143759.72
$ws.Range("H32").Value = 12019.353
$ws.Range("I32").Value = 12066.967
$ws.Range("J32").Value = 11662.25
$ws.Range("K32").Value = 12066.967
$ws.Range("L32").Value = 11662.25
$ws.Range("M32").Value = -11779.967
$ws.Range("N32").Value = -12236.25
$ws.Range("H61").Value = 1834.2903
$ws.Range("I61").Value = 1882.619
$ws.Range("J61").Value = 1732.8
$ws.Range("K61").Value = 1882.619
$ws.Range("L61").Value = 1732.8
$ws.Range("M61").Value = -1670.619
$ws.Range("N61").Value = -2156.8
$ws.Range("H74").Value = 1642.7805
$ws.Range("I74").Value = 1532.5625
$ws.Range("J74").Value = 2034.6666
$ws.Range("K74").Value = 1532.5625
$ws.Range("L74").Value = 2034.6666
$ws.Range("M74").Value = -658.5625
$ws.Range("N74").Value = -3782.6666
$ws.Range("H77").Value = 1642.7805
$ws.Range("I77").Value = 1532.5625
$ws.Range("J77").Value = 2034.6666
$ws.Range("K77").Value = 7662.8125
$ws.Range("L77").Value = 10173.333
$ws.Range("M77").Value = -3294.8125
$ws.Range("N77").Value = -18909.333
$ws.Range("H116").Value = 29345.742
$ws.Range("I116").Value = 798.75
$ws.Range("J116").Value = 143533.72
$ws.Range("K116").Value = 798.75
$ws.Range("L116").Value = 143533.72
$ws.Range("M116").Value = 1495.25
$ws.Range("N116").Value = -148121.72
$ws.Range("H132").Value = 4069.6936
$ws.Range("I132").Value = 4356.636
$ws.Range("J132").Value = 3368.2778
$ws.Range("K132").Value = 13069.908
$ws.Range("L132").Value = 10104.8334
$ws.Range("M132").Value = -10539.908
$ws.Range("N132").Value = -15164.8334
$ws.Range("H136").Value = 1834.2903
$ws.Range("I136").Value = 1882.619
$ws.Range("J136").Value = 1732.8
$ws.Range("K136").Value = 5647.857
$ws.Range("L136").Value = 5198.4
$ws.Range("M136").Value = -3097.857
$ws.Range("N136").Value = -10298.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 29345.742
$ws.Range("I3").Value = 798.75
$ws.Range("J3").Value = 143533.72
$ws.Range("K3").Value = 798.75
$ws.Range("L3").Value = 143533.72
$ws.Range("M3").Value = -684.75
$ws.Range("N3").Value = -143761.72
$ws.Range("H105").Value = 5497457
$ws.Range("I105").Value = 7145603
$ws.Range("J105").Value = 3636.8333
$ws.Range("K105").Value = 7145603
$ws.Range("L105").Value = 3636.8333
$ws.Range("M105").Value = -7143856
$ws.Range("N105").Value = -7130.8333
$ws.Range("H134").Value = 3199.8948
$ws.Range("I134").Value = 2898.6667
$ws.Range("J134").Value = 3716.2856
$ws.Range("K134").Value = 8696.000100000001
$ws.Range("L134").Value = 11148.8568
$ws.Range("M134").Value = -6161.000100000001
$ws.Range("N134").Value = -16218.8568

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2014.0366
$ws.Range("I31").Value = 2186.3333
$ws.Range("J31").Value = 1857.7675
$ws.Range("K31").Value = 2186.3333
$ws.Range("L31").Value = 1857.7675
$ws.Range("M31").Value = -1891.3333
$ws.Range("N31").Value = -2447.7675
$ws.Range("H34").Value = 2014.0366
$ws.Range("I34").Value = 2186.3333
$ws.Range("J34").Value = 1857.7675
$ws.Range("K34").Value = 2186.3333
$ws.Range("L34").Value = 1857.7675
$ws.Range("M34").Value = -1984.3333
$ws.Range("N34").Value = -2261.7675
$ws.Range("H134").Value = 2509.3333
$ws.Range("I134").Value = 2010.1818
$ws.Range("K134").Value = 6030.5454
$ws.Range("M134").Value = -3495.5454

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1379.4932
$ws.Range("I107").Value = 1144.3636
$ws.Range("J107").Value = 1736.2413
$ws.Range("K107").Value = 3433.0908
$ws.Range("L107").Value = 5208.7239
$ws.Range("M107").Value = -1513.0908
$ws.Range("N107").Value = -9048.723900000001
$ws.Range("H122").Value = 628.63635
$ws.Range("I122").Value = 450.05554
$ws.Range("J122").Value = 842.93335
$ws.Range("K122").Value = 4050.49986
$ws.Range("L122").Value = 7586.40015
$ws.Range("M122").Value = -1600.49986
$ws.Range("N122").Value = -12486.40015
$ws.Range("H131").Value = 16669496
$ws.Range("J131").Value = 21742428
$ws.Range("L131").Value = 65227284
$ws.Range("N131").Value = -65237364
$ws.Range("H140").Value = 2019.8276
$ws.Range("J140").Value = 3284.2307
$ws.Range("L140").Value = 9852.6921
$ws.Range("N140").Value = -20212.6921

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5688.4194
$ws.Range("I132").Value = 6016.143
$ws.Range("J132").Value = 5000.2
$ws.Range("K132").Value = 18048.429
$ws.Range("L132").Value = 15000.6
$ws.Range("M132").Value = -15518.429
$ws.Range("N132").Value = -20060.6
$ws.Range("H136").Value = 4262.5
$ws.Range("I136").Value = 2100
$ws.Range("J136").Value = 4571.4287
$ws.Range("K136").Value = 6300
$ws.Range("L136").Value = 13714.2861
$ws.Range("M136").Value = -3750
$ws.Range("N136").Value = -18814.2861

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3236.639
$ws.Range("I132").Value = 3167.4167
$ws.Range("J132").Value = 3375.0833
$ws.Range("K132").Value = 9502.250100000001
$ws.Range("L132").Value = 10125.2499
$ws.Range("M132").Value = -6972.250100000001
$ws.Range("N132").Value = -15185.2499
$ws.Range("H136").Value = 4204.8335
$ws.Range("I136").Value = 5148.375
$ws.Range("J136").Value = 3450
$ws.Range("K136").Value = 15445.125
$ws.Range("L136").Value = 10350
$ws.Range("M136").Value = -12895.125
$ws.Range("N136").Value = -15450

